# FIX months in date fields and ADD drupal from basics to advance with 3 live
# projects course to the course list

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Course")

# --- Fix the month on existing date cells (May -> August) ---
$ws.Range("J2").Value = 45150
$ws.Range("F3").Value = 45148
$ws.Range("J3").Value = 45150
$ws.Range("F4").Value = 45149
$ws.Range("J4").Value = 45149

# --- Add the new course row (row 5) ---
# Columns: A=CourseID, B=Course Name, C=Course Link, D=Course Author,
#          E=Certificate, F=Date Started, G=Date Finished, H=Section,
#          I=Lesson, J=Date
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Drupal From Basics to Advance with 3 Live Projects"
$ws.Range("C5").Value = "https://www.udemy.com/course/advanced-web-development-with-drupal"
$ws.Range("D5").Value = "Abdul Rehman"
$ws.Range("E5").Value = "No"
$ws.Range("F5").Value = 45150
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 45150

$ws.Range("A5:J5").Select()
